$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.400.11"
$ws.Range("E2").Value = "  -1.11%  "
$ws.Range("D3").Value = "2.284.58"
$ws.Range("E3").Value = "  -1.08%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "300.17"
$ws.Range("E5").Value = "  -0.94%  "
$ws.Range("D6").Value = "96.66"
$ws.Range("E6").Value = "  -3.29%  "
$ws.Range("D7").Value = "0.498"
$ws.Range("E7").Value = "  -1.35%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "0.493"
$ws.Range("E9").Value = "  -2.14%  "
$ws.Range("D10").Value = "33.37"
$ws.Range("E10").Value = "  -4.30%  "
$ws.Range("D11").Value = "0.0791"
$ws.Range("E11").Value = "  -0.11%  "
$ws.Range("D12").Value = "48.48"
$ws.Range("E12").Value = "  -6.00%  "
$ws.Range("E13").Value = "  +2.10%  "
$ws.Range("D14").Value = "16.01"
$ws.Range("E14").Value = "  +2.32%  "
$ws.Range("E15").Value = "  -0.17%  "
$ws.Range("D16").Value = "2.635.39"
$ws.Range("E16").Value = "  -1.02%  "
$ws.Range("D17").Value = "2.281.49"
$ws.Range("E17").Value = "  -0.94%  "
$ws.Range("D18").Value = "0.793"
$ws.Range("E18").Value = "  -1.00%  "
$ws.Range("D19").Value = "42.340.35"
$ws.Range("E19").Value = "  -0.92%  "
$ws.Range("D20").Value = "11.68"
$ws.Range("E20").Value = "  +0.31%  "
$ws.Range("D21").Value = "0.0₃0895"
$ws.Range("E21").Value = "  -1.31%  "
$ws.Range("D22").Value = "6.00"
$ws.Range("E22").Value = "  -0.92%  "
$ws.Range("D23").Value = "66.57"
$ws.Range("E23").Value = "  -2.02%  "
$ws.Range("D24").Value = "236.01"
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("E25").Value = "  +0.31%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("E27").Value = "  -2.13%  "
$ws.Range("D28").Value = "23.87"
$ws.Range("E28").Value = "  -3.93%  "
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value = "167.38"
$ws.Range("E29").Value = "  +1.31%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "2.07"
$ws.Range("E30").Value = "  -4.64%  "
$ws.Range("D31").Value = "33.86"
$ws.Range("E31").Value = "  -1.89%  "
$ws.Range("D32").Value = "9.16"
$ws.Range("E32").Value = "  +0.39%  "
$ws.Range("E33").Value = "  +0.20%  "
$ws.Range("D34").Value = "4.68"
$ws.Range("E34").Value = "  +5.05%  "
$ws.Range("D35").Value = "4.95"
$ws.Range("E35").Value = "  -1.41%  "
$ws.Range("D36").Value = "16.91"
$ws.Range("E36").Value = "  +0.55%  "
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").Value = "0.0695"
$ws.Range("E37").Value = "  -0.26%  "
$ws.Range("B38").Value = "WEMIXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D38").Value = "2.34"
$ws.Range("E38").Value = "  -3.11%  "
$ws.Range("D39").Value = "2.80"
$ws.Range("E39").Value = "  -3.20%  "
$ws.Range("D40").Value = "0.0991"
$ws.Range("E40").Value = "  -1.30%  "
$ws.Range("D41").Value = "1.74"
$ws.Range("E41").Value = "  -3.72%  "
$ws.Range("E42").Value = "  -2.13%  "
$ws.Range("D43").Value = "2.27"
$ws.Range("E43").Value = "  -7.94%  "
$ws.Range("D44").Value = "1.960.01"
$ws.Range("E44").Value = "  -0.39%  "
$ws.Range("E45").Value = "  -0.37%  "
$ws.Range("D46").Value = "17.48"
$ws.Range("E46").Value = "  -5.14%  "
$ws.Range("D47").Value = "9.66"
$ws.Range("E47").Value = "  -5.77%  "
$ws.Range("D48").Value = "2.81"
$ws.Range("E48").Value = "  -2.59%  "
$ws.Range("D49").Value = "2.506.03"
$ws.Range("E49").Value = "  -0.98%  "
$ws.Range("D50").Value = "52.43"
$ws.Range("E50").Value = "  -5.67%  "
$ws.Range("E51").Value = "  -3.20%  "
